$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new column H "Is verified" -------------------------------------
# Copy formatting from column G (same style/border/fill) down to column H
# for the full used range (rows 1-50), then set the header + data values.
$ws.Range("G1:G50").Copy()
$ws.Range("H1:H50").PasteSpecial(-4122)

# Match column width of F:G for the new column H.
$ws.Range("H1").ColumnWidth = $ws.Range("G1").ColumnWidth

# Header
$ws.Range("H1").Value = "Is verified"

# Data rows
$ws.Range("H2").Value = "Yes"
$ws.Range("H3").Value = "No"
# H4 intentionally left blank (no value for that contact)
$ws.Range("H5").Value = "yes"
$ws.Range("H6").Value = "no"

# --- 2. Fix the Email value/hyperlink display text in F5 --------------------
$ws.Range("F5").Value = "ttest@example.com"
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 5 -and $h.Range.Column -eq 6) {
        $h.TextToDisplay = "ttest@example.com"
    }
}
